# Update "想去人数" (want-to-go count) values in column F for rows 2, 4, 8, 10, 11, 17, 18, 19
# These updates apply identically to both the "展览" sheet and the "全部类型" sheet,
# since the latter aggregates the rows of the former (plus other sheets).

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 160
    "F4"  = 12298
    "F8"  = 87
    "F10" = 189
    "F11" = 445
    "F17" = 3116
    "F18" = 91
    "F19" = 934
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
